$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E2 and F2 text values
$ws.Range("E2").Value = "반짝이는 미니 크리스마스트리"
$ws.Range("F2").Value = "작은 공간에도 완벽한 45cm 미니트리로 따뜻한 연말 분위기를 만끽하세요."

# Remove column G entirely (G1 header "link" and G2 empty cell)
$ws.Range("G1:G2").Delete()
